# Append one new user-engagement record as row 8 (sheet currently has
# header row 1 + data rows 2-7, i.e. used range A1:K7 -> A1:K8 after this).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 8

$ws.Cells.Item($row, 1).Value  = 1417827147                  # user_id
$ws.Cells.Item($row, 2).Value  = "user_1417827147"            # username
$ws.Cells.Item($row, 3).Value  = -1                           # level
$ws.Cells.Item($row, 4).Value  = ""                           # last_message_date
$ws.Cells.Item($row, 5).Value  = ""                           # last_response
$ws.Cells.Item($row, 6).Value  = "unreachable"                # response_status
$ws.Cells.Item($row, 7).Value  = ""                           # level_3_ai_response
$ws.Cells.Item($row, 8).Value  = $false                       # subscription_checked
$ws.Cells.Item($row, 9).Value  = $false                       # level_4_reminder_sent
$ws.Cells.Item($row, 10).Value = ""                           # decision
$ws.Cells.Item($row, 11).Value = "Added during extraction"    # notes
